$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.910.97'
$ws.Range('D3').Value = '3.736.08'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.53'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('D7').Value = '3.734.86'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.41'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000256'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.89'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '4.362.19'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '3.733.55'
$ws.Range('D17').Value = '67.855.35'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.87'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '463.46'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.694'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.57'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('E25').Value = '  +8.98%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '29.65'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('E33').Value = '  -3.33%  '
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('D36').Value = '3.688.77'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.42'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.994'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.95'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +16.06%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.70'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.33%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.38'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '144.56'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '386.57'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').Value = '2.759.75'
$ws.Range('E51').Value = '  +3.17%  '
